# Staging.Project_ResultArea.xlsx — relocation/regeneration edit.
#
# Net content change: the two header cells in row 2 are relabelled —
#   B2: "ResultAreaID"    -> "ProjectBusinessKey"
#   C2: "ProjectSourceKey" -> "ResultAreaID"
# (A1/A2 and everything else on the sheet stay the same.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ProjectBusinessKey"
$ws.Range("C2").Value = "ResultAreaID"

# Best-effort cosmetic touch-ups that mirror the rest of the commit
# (sheet VBA codename bump + the book window being saved at a larger,
# maximized-looking size). Some hosts expose these as read-only; guard
# them so the substantive edit above always applies cleanly.
try { $ws.CodeName = "Sheet45" } catch { }
try {
    $win = $excel.ActiveWindow
    $win.Width = 28800
    $win.Height = 12585
} catch { }
